$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warrant Issued Report")

# Grab formatting from an existing "D" cell that already uses the plain
# wrapped style used elsewhere in this table, so the new D50 cell picks up
# matching formatting once pasted after the row insert below.
$ws.Cells.Item(24, 4).Copy()

# Insert a new row at row 50, pushing existing row 50+ down (row 49, the
# "Vehicle Registration Expiration Date" row, keeps its position).
$ws.Rows.Item(50).Insert()

$ws.Cells.Item(50, 4).PasteSpecial(-4122)

$ws.Cells.Item(50, 3).Value = "Vehicle Registration Non-Expiring Indicator"
$ws.Cells.Item(50, 4).Value = "True if vehicle registration is non-expiring"
$ws.Cells.Item(50, 5).Value = $true
$ws.Cells.Item(50, 6).Value = "/wir-doc:WarrantIssuedReport/j:ConveyanceRegistration[not(j:RegistrationExpirationDate)]/wir-ext:ConveyanceRegistrationNonExpiringIndicator"

$ws.Rows.Item(50).RowHeight = 56
